$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Person")
$ws.Range("A2").Value = "UKT_0001"
